$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = New-Object 'object[,]' 57,3
$values[0,0] = 62542; $values[0,1] = 8974; $values[0,2] = 10094
$values[1,0] = 29626; $values[1,1] = 3380; $values[1,2] = 3526
$values[2,0] = 104519; $values[2,1] = 7355; $values[2,2] = 9564
$values[3,0] = 2207; $values[3,1] = 677; $values[3,2] = 113
$values[4,0] = 61415; $values[4,1] = 10611; $values[4,2] = 9298
$values[5,0] = 7557; $values[5,1] = 1583; $values[5,2] = 1539
$values[6,0] = 7950; $values[6,1] = 1388; $values[6,2] = 1028
$values[7,0] = 2983; $values[7,1] = 494; $values[7,2] = 134
$values[8,0] = 435; $values[8,1] = 359; $values[8,2] = 8
$values[9,0] = 0; $values[9,1] = 0; $values[9,2] = 0
$values[10,0] = 1346; $values[10,1] = 402; $values[10,2] = 279
$values[11,0] = 4218; $values[11,1] = 1841; $values[11,2] = 1296
$values[12,0] = 8403; $values[12,1] = 2850; $values[12,2] = 1252
$values[13,0] = 4654; $values[13,1] = 2014; $values[13,2] = 810
$values[14,0] = 2747; $values[14,1] = 918; $values[14,2] = 224
$values[15,0] = 23953; $values[15,1] = 3688; $values[15,2] = 4339
$values[16,0] = 4102; $values[16,1] = 1121; $values[16,2] = 799
$values[17,0] = 29470; $values[17,1] = 3222; $values[17,2] = 5050
$values[18,0] = 608; $values[18,1] = 464; $values[18,2] = 29
$values[19,0] = 24844; $values[19,1] = 2950; $values[19,2] = 3893
$values[20,0] = 1699; $values[20,1] = 674; $values[20,2] = 266
$values[21,0] = 25513; $values[21,1] = 4290; $values[21,2] = 4446
$values[22,0] = 98168; $values[22,1] = 11012; $values[22,2] = 10848
$values[23,0] = 7128; $values[23,1] = 2570; $values[23,2] = 981
$values[24,0] = 0; $values[24,1] = 0; $values[24,2] = 0
$values[25,0] = 7998; $values[25,1] = 1553; $values[25,2] = 1514
$values[26,0] = 3089; $values[26,1] = 478; $values[26,2] = 642
$values[27,0] = 22182; $values[27,1] = 4040; $values[27,2] = 3947
$values[28,0] = 808; $values[28,1] = 291; $values[28,2] = 320
$values[29,0] = 3425; $values[29,1] = 2192; $values[29,2] = 446
$values[30,0] = 21145; $values[30,1] = 4289; $values[30,2] = 3813
$values[31,0] = 15459; $values[31,1] = 3622; $values[31,2] = 3723
$values[32,0] = 8881; $values[32,1] = 1073; $values[32,2] = 1878
$values[33,0] = 74093; $values[33,1] = 8088; $values[33,2] = 7612
$values[34,0] = 11873; $values[34,1] = 3697; $values[34,2] = 1952
$values[35,0] = 28118; $values[35,1] = 2720; $values[35,2] = 4068
$values[36,0] = 1260; $values[36,1] = 1078; $values[36,2] = 205
$values[37,0] = 2176; $values[37,1] = 343; $values[37,2] = 1001
$values[38,0] = 3328; $values[38,1] = 425; $values[38,2] = 211
$values[39,0] = 13390; $values[39,1] = 302; $values[39,2] = 325
$values[40,0] = 377; $values[40,1] = 132; $values[40,2] = 85
$values[41,0] = 1124; $values[41,1] = 60; $values[41,2] = 32
$values[42,0] = 2809; $values[42,1] = 195; $values[42,2] = 73
$values[43,0] = 4569; $values[43,1] = 1306; $values[43,2] = 716
$values[44,0] = 17447; $values[44,1] = 4523; $values[44,2] = 3320
$values[45,0] = 41847; $values[45,1] = 4410; $values[45,2] = 6262
$values[46,0] = 20627; $values[46,1] = 4637; $values[46,2] = 1695
$values[47,0] = 15117; $values[47,1] = 1463; $values[47,2] = 2588
$values[48,0] = 39547; $values[48,1] = 4047; $values[48,2] = 4873
$values[49,0] = 5467; $values[49,1] = 643; $values[49,2] = 1239
$values[50,0] = 19291; $values[50,1] = 3983; $values[50,2] = 3221
$values[51,0] = 3081; $values[51,1] = 1026; $values[51,2] = 1539
$values[52,0] = 2938; $values[52,1] = 1903; $values[52,2] = 282
$values[53,0] = 4585; $values[53,1] = 1435; $values[53,2] = 1470
$values[54,0] = 18621; $values[54,1] = 7107; $values[54,2] = 3637
$values[55,0] = 21028; $values[55,1] = 1114; $values[55,2] = 726
$values[56,0] = 912138; $values[56,1] = 136791; $values[56,2] = 128868

$range = $ws.Range("B3:D59")
$range.Value = $values
